$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column I: "Other found locations"
$ws.Cells.Item(1, 9).Value = 'Other found locations'

# Update Authors (column E) text for rows whose author list was reprocessed
$ws.Cells.Item(2, 5).Value = '[ Thanos%Karatzias%null%2,     Elsie%Yan%null%1,     Sally%Jowett%null%1,  Thanos%Karatzias%null%0,  Elsie%Yan%null%1,  Sally%Jowett%null%1]'
$ws.Cells.Item(4, 5).Value = '[Alain%Vanasse%NULL%1,    Alan%Cohen%NULL%1,    Josiane%Courteau%NULL%1,    Patrick%Bergeron%NULL%1,    Roxanne%Dault%NULL%1,    Pierre%Gosselin%NULL%1,    Claudia%Blais%NULL%1,    Diane%Bélanger%NULL%1,    Louis%Rochette%NULL%1,    Fateh%Chebana%NULL%1,    Jan%Semenza%NULL%2,    Jan%Semenza%NULL%0]'
$ws.Cells.Item(5, 5).Value = '[ Leonardo%Trasande%null%2,     Tony T.%Koshy%null%1,     Joseph%Gilbert%null%1,     Lauren K.%Burdine%null%1,     Michael%Marmor%null%1,     Xiaoxia%Han%null%1,     Yongzhao%Shao%null%1,     Claude%Chemtob%null%1,     Teresa M.%Attina%null%1,     Elaine M.%Urbina%null%1,  Leonardo%Trasande%null%0,  Tony T.%Koshy%null%1,  Joseph%Gilbert%null%1,  Lauren K.%Burdine%null%1,  Michael%Marmor%null%1,  Xiaoxia%Han%null%1,  Yongzhao%Shao%null%1,  Claude%Chemtob%null%1,  Teresa M.%Attina%null%1,  Elaine M.%Urbina%null%1]'
$ws.Cells.Item(6, 5).Value = '[Ninon A.%Becquart%NULL%1,    Elena N.%Naumova%NULL%1,    Gitanjali%Singh%NULL%1,    Kenneth K. H.%Chui%NULL%1]'
$ws.Cells.Item(7, 5).Value = '[Hannah T.%Jordan%NULL%1,    Steven D.%Stellman%NULL%1,    Alfredo%Morabia%NULL%1,    Sara A.%Miller‐Archie%NULL%1,    Howard%Alper%NULL%1,    Zoey%Laskaris%NULL%1,    Robert M.%Brackbill%NULL%1,    James E.%Cone%NULL%1]'
$ws.Cells.Item(9, 5).Value = '[Hector%F. Myers%coreGivesNoEmail%1,   Karen%A. Matthews%coreGivesNoEmail%1,   Nataria%T. Joseph%coreGivesNoEmail%1]'
$ws.Cells.Item(10, 5).Value = '[ Nathan%McKinney%null%2,     Chris%Houser%null%1,     Klaus%Meyer-Arendt%null%1,  Nathan%McKinney%null%0,  Chris%Houser%null%1,  Klaus%Meyer-Arendt%null%1]'
$ws.Cells.Item(12, 5).Value = '[ Zhen%Jiao%null%2,     Socrates V.%Kakoulides%null%1,     John%Moscona%null%1,     Jabar%Whittier%null%1,     Sudesh%Srivastav%null%1,     Patrice%Delafontaine%null%2,     Anand%Irimpen%null%1,  Zhen%Jiao%null%0,  Socrates V.%Kakoulides%null%1,  John%Moscona%null%1,  Jabar%Whittier%null%1,  Sudesh%Srivastav%null%1,  Patrice%Delafontaine%null%2,  Anand%Irimpen%null%1]'
$ws.Cells.Item(14, 5).Value = '[ Yasuhiro%Nagayoshi%null%2,     Shinya%Yumoto%null%1,     Kazuhisa%Sakaguchi%null%1,     Chiharu%Shudo%null%1,     Shiro%Takino%null%1,     Motohiro%Hashiyama%null%1,     Yutaka%Kai%null%1,     Yutaka%Kuroda%null%1,     Hiroaki%Kawano%null%1,     Hisao%Ogawa%null%2,  Yasuhiro%Nagayoshi%null%0,  Shinya%Yumoto%null%1,  Kazuhisa%Sakaguchi%null%1,  Chiharu%Shudo%null%1,  Shiro%Takino%null%1,  Motohiro%Hashiyama%null%1,  Yutaka%Kai%null%1,  Yutaka%Kuroda%null%1,  Hiroaki%Kawano%null%1,  Hisao%Ogawa%null%1]'
$ws.Cells.Item(15, 5).Value = '[ Hannah T.%Jordan%null%2,     Sara A.%Miller-Archie%null%1,     James E.%Cone%null%1,     Alfredo%Morabia%null%1,     Steven D.%Stellman%null%1,  Hannah T.%Jordan%null%0,  Sara A.%Miller-Archie%null%1,  James E.%Cone%null%1,  Alfredo%Morabia%null%1,  Steven D.%Stellman%null%1]'
$ws.Cells.Item(17, 5).Value = '[Vivian A.%Fonseca%NULL%1,    Hayden%Smith%NULL%1,    Nitesh%Kuhadiya%NULL%1,    Sharice M.%Leger%NULL%1,    C. Lillian%Yau%NULL%1,    Kristi%Reynolds%NULL%1,    Lizheng%Shi%NULL%1,    Roberta H.%McDuffie%NULL%1,    Tina%Thethi%NULL%1,    Jennifer%John-Kalarickal%NULL%1]'
$ws.Cells.Item(18, 5).Value = '[ Federico%Silva-Palacios%null%2,     Ana Isabel%Casanegra%null%1,     Alan%Shapiro%null%1,     Minh%Phan%null%1,     Beau%Hawkins%null%1,     Ji%Li%null%1,     Julie%Stoner%null%1,     Alfonso%Tafur%null%1,  Federico%Silva-Palacios%null%0,  Ana Isabel%Casanegra%null%1,  Alan%Shapiro%null%1,  Minh%Phan%null%1,  Beau%Hawkins%null%1,  Ji%Li%null%1,  Julie%Stoner%null%1,  Alfonso%Tafur%null%1]'
$ws.Cells.Item(19, 5).Value = '[Tran Huu%Bich%NULL%1,    La Ngoc%Quang%NULL%1,    Le Thi%Thanh Ha%NULL%1,    Tran Thi%Duc Hanh%NULL%1,    Debarati%Guha-Sapir%NULL%1]'
$ws.Cells.Item(20, 5).Value = '[Bijlsma%J A%coreGivesNoEmail%1,   Bouter%L M%coreGivesNoEmail%1,   Bramsen%I%coreGivesNoEmail%1,   Huizink%A C%coreGivesNoEmail%1,   Slottje%P%coreGivesNoEmail%1,   Smid%T%coreGivesNoEmail%1,   Smidt%N%coreGivesNoEmail%1,   Twisk%J W R%coreGivesNoEmail%1,   van%der Ploeg H M%coreGivesNoEmail%1,   van%Mechelen W%coreGivesNoEmail%1,   Witteveen%A B%coreGivesNoEmail%1]'
$ws.Cells.Item(21, 5).Value = '[ Hannah T%Jordan%null%2,     Robert M%Brackbill%null%1,     James E%Cone%null%1,     Indira%Debchoudhury%null%1,     Mark R%Farfel%null%1,     Carolyn M%Greene%null%1,     James L%Hadler%null%1,     Joseph%Kennedy%null%1,     Jiehui%Li%null%1,     Jonathan%Liff%null%1,     Leslie%Stayner%null%1,     Steven D%Stellman%null%1,  Hannah T%Jordan%null%0,  Robert M%Brackbill%null%1,  James E%Cone%null%1,  Indira%Debchoudhury%null%1,  Mark R%Farfel%null%1,  Carolyn M%Greene%null%1,  James L%Hadler%null%1,  Joseph%Kennedy%null%1,  Jiehui%Li%null%1,  Jonathan%Liff%null%1,  Leslie%Stayner%null%1,  Steven D%Stellman%null%1]'
$ws.Cells.Item(23, 5).Value = '[ Matthew N.%Peters%null%2,     Morgan J.%Katz%null%1,     John C.%Moscona%null%1,     Mohi E.%Alkadri%null%1,     Rashad H.%Khazi Syed%null%1,     Thomas A.%Turnage%null%1,     Vikram S.%Nijjar%null%1,     Mohannad B.%Bisharat%null%1,     Patrice%Delafontaine%null%0,     Anand M.%Irimpen%null%1,  Matthew N.%Peters%null%0,  Morgan J.%Katz%null%1,  John C.%Moscona%null%1,  Mohi E.%Alkadri%null%1,  Rashad H.%Khazi Syed%null%1,  Thomas A.%Turnage%null%1,  Vikram S.%Nijjar%null%1,  Mohannad B.%Bisharat%null%1,  Patrice%Delafontaine%null%0,  Anand M.%Irimpen%null%1]'
$ws.Cells.Item(24, 5).Value = '[Ibrahim Baimba%Koroma%NULL%1,    Dena%Javadi%NULL%2,    Dena%Javadi%NULL%0,    Katrina%Hann%NULL%5,    Anthony D%Harries%NULL%2,    Anthony D%Harries%NULL%0,    Francis%Smart%NULL%2,    Francis%Smart%NULL%0,    Thomas%Samba%NULL%1,    Jeffery%Edwards%NULL%4,    Jeffery%Edwards%NULL%0,    Priyakanta%Nayak%NULL%2,    Priyakanta%Nayak%NULL%0,    Katrina%Hann%NULL%0,    Katrina%Hann%NULL%0,    Jeffery%Edwards%NULL%0,    Jeffery%Edwards%NULL%0,    Katrina%Hann%NULL%0,    Katrina%Hann%NULL%0]'
$ws.Cells.Item(29, 5).Value = '[ M. Natalie%Husarewycz%null%2,     Renée%El-Gabalawy%null%1,     Sarvesh%Logsetty%null%1,     Jitender%Sareen%null%1,  M. Natalie%Husarewycz%null%0,  Renée%El-Gabalawy%null%1,  Sarvesh%Logsetty%null%1,  Jitender%Sareen%null%1]'
$ws.Cells.Item(31, 5).Value = '[ Grégoire%Rey%null%2,     Eric%Jougla%null%1,     Anne%Fouillet%null%1,     Gérard%Pavillon%null%1,     Pierre%Bessemoulin%null%1,     Philippe%Frayssinet%null%1,     Jacqueline%Clavel%null%1,     Denis%Hémon%null%1,  Grégoire%Rey%null%0,  Eric%Jougla%null%1,  Anne%Fouillet%null%1,  Gérard%Pavillon%null%1,  Pierre%Bessemoulin%null%1,  Philippe%Frayssinet%null%1,  Jacqueline%Clavel%null%1,  Denis%Hémon%null%1]'
$ws.Cells.Item(32, 5).Value = '[Ahmad%Jabbar%coreGivesNoEmail%1,   Alaa%Boulad%coreGivesNoEmail%1,   Anand%Irimpen%coreGivesNoEmail%1,   Bradley%Deere%coreGivesNoEmail%1,   Christopher%Westley%coreGivesNoEmail%1,   Hassan%Baydoun%coreGivesNoEmail%1,   Holly%Gonzales%coreGivesNoEmail%1,   Indrajeet%Mahata%coreGivesNoEmail%1,   John%C. Moscona%coreGivesNoEmail%1,   Kapil%Yadav%coreGivesNoEmail%1,   Matthew%N. Peters%coreGivesNoEmail%1,   Patrick%Ters%coreGivesNoEmail%1,   Paul%Katigbak%coreGivesNoEmail%1,   Rohit%Maini%coreGivesNoEmail%1,   Ryan%Nelson%coreGivesNoEmail%1,   Sudesh%Srivastav%coreGivesNoEmail%1,   Taraka%V. Gadiraju%coreGivesNoEmail%1]'

# Rows whose Doi could not be matched to an ID any more -> ID/ID Format reset to not found/N/A
$ws.Cells.Item(9, 6).Value = "not found"
$ws.Cells.Item(9, 7).Value = "N/A"
$ws.Cells.Item(20, 6).Value = "not found"
$ws.Cells.Item(20, 7).Value = "N/A"
$ws.Cells.Item(30, 6).Value = "not found"
$ws.Cells.Item(30, 7).Value = "N/A"
$ws.Cells.Item(32, 6).Value = "not found"
$ws.Cells.Item(32, 7).Value = "N/A"

# New column I ("Other found locations") values per row; "_PMC" flags rows whose record was located via PMC, blank otherwise
$ws.Cells.Item(2, 9).Value = ""
$ws.Cells.Item(3, 9).Value = ""
$ws.Cells.Item(4, 9).Value = "_PMC"
$ws.Cells.Item(5, 9).Value = ""
$ws.Cells.Item(6, 9).Value = "_PMC"
$ws.Cells.Item(7, 9).Value = "_PMC"
$ws.Cells.Item(8, 9).Value = ""
$ws.Cells.Item(9, 9).Value = ""
$ws.Cells.Item(10, 9).Value = ""
$ws.Cells.Item(11, 9).Value = ""
$ws.Cells.Item(12, 9).Value = ""
$ws.Cells.Item(13, 9).Value = ""
$ws.Cells.Item(14, 9).Value = ""
$ws.Cells.Item(15, 9).Value = ""
$ws.Cells.Item(16, 9).Value = ""
$ws.Cells.Item(17, 9).Value = "_PMC"
$ws.Cells.Item(18, 9).Value = ""
$ws.Cells.Item(19, 9).Value = "_PMC"
$ws.Cells.Item(20, 9).Value = ""
$ws.Cells.Item(21, 9).Value = ""
$ws.Cells.Item(22, 9).Value = ""
$ws.Cells.Item(23, 9).Value = ""
$ws.Cells.Item(24, 9).Value = "_PMC"
$ws.Cells.Item(25, 9).Value = ""
$ws.Cells.Item(26, 9).Value = ""
$ws.Cells.Item(27, 9).Value = ""
$ws.Cells.Item(28, 9).Value = ""
$ws.Cells.Item(29, 9).Value = ""
$ws.Cells.Item(30, 9).Value = ""
$ws.Cells.Item(31, 9).Value = ""
$ws.Cells.Item(32, 9).Value = ""
$ws.Cells.Item(33, 9).Value = ""
